# Update Gilberto's schedule: shift some class entries to their new time slots/days.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B6").Value  = "MEC-2B-Mec. Tec. Res. Mat."
$ws.Range("C7").Value  = "-"

$ws.Range("B14").Value = "-"

$ws.Range("B15").Value = "-"
$ws.Range("E15").Value = "MEC-2A-Mec. Tec. Res. Mat."

$ws.Range("E16").Value = "MEC-2A-Mec. Tec. Res. Mat."

$ws.Range("C18").Value = "-"
$ws.Range("C19").Value = "MEC-1NA-M.T.R.M."

$ws.Range("C20").Value = "-"
$ws.Range("E20").Value = "-"
$ws.Range("F20").Value = "MEC-1NB-M.T.R.M."

$ws.Range("B21").Value = "ELM-1NA-Tecnologias Mecânicas"
$ws.Range("C21").Value = "-"
$ws.Range("F21").Value = "ELM-1NA-Tecnologias Mecânicas"
